$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.968.94"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.779.52"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'226.01"
$ws.Range("E5").Value = "  +2.65%  "
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "'32.21"
$ws.Range("E8").Value = "  +2.82%  "
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("D10").Value = "'0.0703"
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("D11").Value = "'0.0936"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").Value = "2.034.12"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.788.79"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'10.91"
$ws.Range("E14").Value = "  +3.57%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "33.928.35"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.619"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "'242.83"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "0.0₃0782"
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "'10.66"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  -2.70%  "
$ws.Range("D25").Value = "'160.16"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("D26").Value = "'16.27"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").Value = "'7.06"
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").Value = "'0.113"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  +3.72%  "
$ws.Range("D31").Value = "'0.0512"
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("D32").Value = "'3.64"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("D34").Value = "'1.81"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("D35").Value = "1.392.23"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").Value = "'0.653"
$ws.Range("E36").Value = "  +4.83%  "
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("E38").Value = "  +0.90%  "
$ws.Range("D39").Value = "'2.36"
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").Value = "'2.20"
$ws.Range("E40").Value = "  +5.10%  "
$ws.Range("D41").Value = "'0.909"
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("D42").Value = "'2.67"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").Value = "'77.46"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").Value = "'13.14"
$ws.Range("E44").Value = "  +13.26%  "
$ws.Range("E45").Value = "  +2.77%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'107.77"
$ws.Range("E46").Value = "  +2.57%  "
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").Value = "'0.0495"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("E48").Value = "  +11.63%  "
$ws.Range("D49").Value = "'5.81"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").Value = "1.933.70"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("E51").Value = "  +0.45%  "
